$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 48383
$ws.Range("B2").Value = "Yasmin Lopes"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45099
$ws.Range("G2").Value = 10645.12

# Row 3
$ws.Range("A3").Value = 30772
$ws.Range("B3").Value = "Sabrina da Cruz"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Viagem de negócios"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45102
$ws.Range("G3").Value = 3350.43

# Row 4
$ws.Range("A4").Value = 2931
$ws.Range("B4").Value = "Davi Lucca Fogaça"
$ws.Range("C4").Value = "Recursos Humanos"
$ws.Range("D4").Value = "Viagem de negócios"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45084
$ws.Range("G4").Value = 3822.63

# Row 5
$ws.Range("A5").Value = 58231
$ws.Range("B5").Value = "Dr. Paulo da Cruz"
$ws.Range("C5").Value = "Engenharia"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45080
$ws.Range("G5").Value = 8303.46

# Row 6
$ws.Range("A6").Value = 23332
$ws.Range("B6").Value = "Lívia da Mota"
$ws.Range("C6").Value = "P&D"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45098
$ws.Range("G6").Value = 5424.79

# Row 7
$ws.Range("A7").Value = 21625
$ws.Range("B7").Value = "Daniela Gonçalves"
$ws.Range("C7").Value = "TI"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45105
$ws.Range("G7").Value = 11925.76

# Row 8
$ws.Range("A8").Value = 30455
$ws.Range("B8").Value = "Sra. Larissa Duarte"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45101
$ws.Range("G8").Value = 12383.5

# Row 9
$ws.Range("A9").Value = 63032
$ws.Range("B9").Value = "Cauê Pires"
$ws.Range("C9").Value = "TI"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45093
$ws.Range("G9").Value = 3509.08

# Row 10
$ws.Range("A10").Value = 3349
$ws.Range("B10").Value = "Brenda Fernandes"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Consulta médica"
$ws.Range("F10").Value = 45104
$ws.Range("G10").Value = 11312.91

# Row 11
$ws.Range("A11").Value = 71162
$ws.Range("B11").Value = "Leandro Duarte"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45078
$ws.Range("G11").Value = 7057.89
